$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column E with a date header (serial 45354 = 2024-03-03) and values
$ws.Range("E1").Value = 45354
$ws.Range("E1").NumberFormat = "mm-dd-yy"

$ws.Range("E2").Value = 3545
$ws.Range("E3").Value = 2975
$ws.Range("E4").Value = 4050
$ws.Range("E5").Value = 7118
$ws.Range("E6").Value = 216

# Size column E to fit the new date column (mirrors Excel's AutoFit result)
$ws.Columns.Item(5).ColumnWidth = 9.43

# Update selection to match the new active cell
$ws.Range("E6").Select()
